# Auto-generated Excel COM-interop edit script
# Applies targeted value updates/additions/deletions to H:N columns
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR to match the
# refreshed market-board price snapshot described in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4476.8
$ws.Range("J19").Value = 4472.6665
$ws.Range("L19").Value = 4472.6665
$ws.Range("N19").Value = -4822.6665

$ws.Range("H32").Value = 5111.1763
$ws.Range("J32").Value = 4545.9287
$ws.Range("L32").Value = 4545.9287
$ws.Range("N32").Value = -5197.9287

$ws.Range("H51").Value = 9049.549999999999
$ws.Range("J51").Value = 9312
$ws.Range("L51").Value = 9312
$ws.Range("N51").Value = -10280

$ws.Range("H53").Value = 634.85
$ws.Range("I53").Value = 669.5454999999999
$ws.Range("K53").Value = 669.5454999999999
$ws.Range("M53").Value = -32.54549999999995

$ws.Range("H69").Value = 6986
$ws.Range("J69").Value = 6986
$ws.Range("L69").Value = 20958
$ws.Range("N69").Value = -22706

$ws.Range("H72").Value = 6986
$ws.Range("J72").Value = 6986
$ws.Range("L72").Value = 62874
$ws.Range("N72").Value = -71610

$ws.Range("H92").Value = 201.81818
$ws.Range("I92").Value = 201.81818
$ws.Range("K92").Value = 201.81818
$ws.Range("M92").Value = 1046.18182

$ws.Range("H112").Value = 3094.5715
$ws.Range("J112").Value = 3279.4167
$ws.Range("L112").Value = 9838.250100000001
$ws.Range("N112").Value = -12054.2501

$ws.Range("H113").Value = 142861420
$ws.Range("J113").Value = 8699.5
$ws.Range("L113").Value = 8699.5
$ws.Range("N113").Value = -15207.5

$ws.Range("H127").Value = 1002
$ws.Range("I127").Value = 1002
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 3006
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 1954
$ws.Range("N127").ClearContents()

$ws.Range("H135").Value = 402.34616
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H138").Value = 66669264
$ws.Range("I138").Value = 2867.0908
$ws.Range("J138").Value = 250001860
$ws.Range("K138").Value = 8601.2724
$ws.Range("L138").Value = 750005580
$ws.Range("M138").Value = -3461.2724
$ws.Range("N138").Value = -750015860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3966.5386
$ws.Range("I2").Value = 1001.2
$ws.Range("J2").Value = 5819.875
$ws.Range("K2").Value = 1001.2
$ws.Range("L2").Value = 5819.875
$ws.Range("M2").Value = -888.2
$ws.Range("N2").Value = -6045.875

$ws.Range("H45").Value = 6823.1816
$ws.Range("I45").Value = 4011
$ws.Range("J45").Value = 9166.666999999999
$ws.Range("K45").Value = 4011
$ws.Range("L45").Value = 9166.666999999999
$ws.Range("M45").Value = -3634
$ws.Range("N45").Value = -9920.666999999999

$ws.Range("H55").Value = 24737.5
$ws.Range("J55").Value = 34500
$ws.Range("L55").Value = 34500
$ws.Range("N55").Value = -35130

$ws.Range("H97").Value = 984.55554
$ws.Range("I97").Value = 951.3333
$ws.Range("K97").Value = 951.3333
$ws.Range("M97").Value = -455.3333

$ws.Range("H116").Value = 3966.5386
$ws.Range("I116").Value = 1001.2
$ws.Range("J116").Value = 5819.875
$ws.Range("K116").Value = 1001.2
$ws.Range("L116").Value = 5819.875
$ws.Range("M116").Value = 1292.8
$ws.Range("N116").Value = -10407.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3966.5386
$ws.Range("I3").Value = 1001.2
$ws.Range("J3").Value = 5819.875
$ws.Range("K3").Value = 1001.2
$ws.Range("L3").Value = 5819.875
$ws.Range("M3").Value = -887.2
$ws.Range("N3").Value = -6047.875

$ws.Range("H20").Value = 5438.8667
$ws.Range("I20").Value = 4620.5557
$ws.Range("K20").Value = 4620.5557
$ws.Range("M20").Value = -4373.5557

$ws.Range("H80").Value = 903.2
$ws.Range("J80").Value = 1004
$ws.Range("L80").Value = 1004
$ws.Range("N80").Value = -3000

$ws.Range("H83").Value = 903.2
$ws.Range("J83").Value = 1004
$ws.Range("L83").Value = 5020
$ws.Range("N83").Value = -15004

$ws.Range("H107").Value = 7395.727
$ws.Range("I107").Value = 7385.5
$ws.Range("K107").Value = 7385.5
$ws.Range("M107").Value = -5465.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 40992.25
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

$ws.Range("H62").Value = 7166.6665
$ws.Range("I62").Value = 6750
$ws.Range("K62").Value = 6750
$ws.Range("M62").Value = -6126

$ws.Range("H65").Value = 7166.6665
$ws.Range("I65").Value = 6750
$ws.Range("K65").Value = 33750
$ws.Range("M65").Value = -30630

$ws.Range("H68").Value = 41271.57
$ws.Range("J68").Value = 41271.57
$ws.Range("L68").Value = 41271.57
$ws.Range("N68").Value = -42769.57

$ws.Range("H71").Value = 41271.57
$ws.Range("J71").Value = 41271.57
$ws.Range("L71").Value = 123814.71
$ws.Range("N71").Value = -131302.71

$ws.Range("H86").Value = 11601.857
$ws.Range("I86").Value = 11843.2
$ws.Range("J86").Value = 10998.5
$ws.Range("K86").Value = 11843.2
$ws.Range("L86").Value = 10998.5
$ws.Range("M86").Value = -10720.2
$ws.Range("N86").Value = -13244.5

$ws.Range("H89").Value = 11601.857
$ws.Range("I89").Value = 11843.2
$ws.Range("J89").Value = 10998.5
$ws.Range("K89").Value = 59216
$ws.Range("L89").Value = 54992.5
$ws.Range("M89").Value = -53600
$ws.Range("N89").Value = -66224.5

$ws.Range("H107").Value = 456.97058
$ws.Range("J107").Value = 541
$ws.Range("L107").Value = 541
$ws.Range("N107").Value = -4381

$ws.Range("H134").Value = 6163.4287
$ws.Range("I134").Value = 5799.8184
$ws.Range("K134").Value = 17399.4552
$ws.Range("M134").Value = -14864.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 379.5
$ws.Range("J29").Value = 194
$ws.Range("L29").Value = 582
$ws.Range("N29").Value = -1136

$ws.Range("H128").Value = 516698.2
$ws.Range("I128").Value = 516698.2
$ws.Range("K128").Value = 1550094.6
$ws.Range("M128").Value = -1545114.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 363.33334
$ws.Range("I17").Value = 363.33334
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 363.33334
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -195.33334
$ws.Range("N17").ClearContents()

$ws.Range("H113").Value = 3068.1428
$ws.Range("I113").Value = 3497
$ws.Range("J113").Value = 2496.3333
$ws.Range("K113").Value = 3497
$ws.Range("L113").Value = 2496.3333
$ws.Range("M113").Value = -1327
$ws.Range("N113").Value = -6836.3333

$ws.Range("H126").Value = 5250
$ws.Range("I126").Value = 3750
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 11250
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -8780
$ws.Range("N126").Value = -22940

$ws.Range("H132").Value = 3302.0356
$ws.Range("I132").Value = 2844.682
$ws.Range("K132").Value = 8534.045999999998
$ws.Range("M132").Value = -6004.045999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3010.9
$ws.Range("I40").Value = 2443.5715
$ws.Range("K40").Value = 2443.5715
$ws.Range("M40").Value = -2307.5715

$ws.Range("H46").Value = 7345.8335
$ws.Range("I46").Value = 3551.4614
$ws.Range("K46").Value = 3551.4614
$ws.Range("M46").Value = -3363.4614

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H68").Value = 2125
$ws.Range("I68").Value = 2150
$ws.Range("K68").Value = 2150
$ws.Range("M68").Value = -1401

$ws.Range("H71").Value = 2125
$ws.Range("I71").Value = 2150
$ws.Range("K71").Value = 10750
$ws.Range("M71").Value = -7006

$ws.Range("H122").Value = 3399.6843
$ws.Range("I122").Value = 3338.5557
$ws.Range("K122").Value = 10015.6671
$ws.Range("M122").Value = -7565.667099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 28000
$ws.Range("J40").Value = 28000
$ws.Range("L40").Value = 28000
$ws.Range("N40").Value = -28298

$ws.Range("H41").Value = 11931.214
$ws.Range("I41").Value = 9000
$ws.Range("J41").Value = 12156.692
$ws.Range("K41").Value = 9000
$ws.Range("L41").Value = 12156.692
$ws.Range("M41").Value = -8610
$ws.Range("N41").Value = -12936.692

$ws.Range("H107").Value = 601.7895
$ws.Range("I107").Value = 537.8
$ws.Range("J107").Value = 841.75
$ws.Range("K107").Value = 1613.4
$ws.Range("L107").Value = 2525.25
$ws.Range("M107").Value = 306.6000000000001
$ws.Range("N107").Value = -6365.25

$ws.Range("H113").Value = 1118.3334
$ws.Range("I113").Value = 1218.5454
$ws.Range("K113").Value = 3655.6362
$ws.Range("M113").Value = -1485.6362

$ws.Range("H116").Value = 199975
$ws.Range("J116").Value = 199975
$ws.Range("L116").Value = 199975
$ws.Range("N116").Value = -209153

$ws.Range("H122").Value = 3437.7
$ws.Range("I122").Value = 1547.125
$ws.Range("J122").Value = 11000
$ws.Range("K122").Value = 4641.375
$ws.Range("L122").Value = 33000
$ws.Range("M122").Value = -2191.375
$ws.Range("N122").Value = -37900

$ws.Range("H132").Value = 3385.56
$ws.Range("I132").Value = 3562
$ws.Range("K132").Value = 10686
$ws.Range("M132").Value = -8156
